$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-10-04 Wednesday" "2023-10-05 Thursday"

Replace-Text "24×56=" "11×12="
Replace-Text "94×80=" "59×98="
Replace-Text "21×49=" "56×34="
Replace-Text "52×19=" "99×85="
Replace-Text "76×19=" "92×44="
Replace-Text "69×56=" "53×78="
Replace-Text "59×39=" "65×68="
Replace-Text "70×39=" "66×86="
Replace-Text "93×44=" "17×47="
Replace-Text "70×21=" "38×21="
Replace-Text "13×47=" "38×12="
Replace-Text "79×69=" "31×64="
Replace-Text "98×60=" "30×91="
Replace-Text "53×29=" "84×32="
Replace-Text "27×97=" "28×11="
Replace-Text "39×47=" "16×38="
Replace-Text "76×14=" "18×90="
Replace-Text "37×21=" "25×16="
Replace-Text "42×29=" "85×40="
Replace-Text "66×76=" "37×11="
Replace-Text "94×64=" "34×46="
Replace-Text "43×20=" "35×47="
Replace-Text "30×25=" "42×97="
Replace-Text "46×21=" "11×34="
Replace-Text "32×43=" "62×74="
